$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (date moved from 03-20 to 03-21)
$ws.Name = "Through 2022-03-21"

# Update header label in I1 (shared string text)
$ws.Range("I1").Value = "2022 (through 03-21)"

# Update data cells per diff
$ws.Range("I4").Value = 90     # March 2022 total
$ws.Range("H6").Value = 108    # May 2021
$ws.Range("H14").Value = 1852  # 2021 Total
$ws.Range("I14").Value = 390   # 2022 Total
